$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shared-string pool indices are assigned in first-use order, so write the
# new cell values in the same order the target workbook's sharedStrings.xml
# lists them: B7, B9, B10, B6, B8.
$ws.Range("B7").Value = "Verify that Paperback/Hardcover/kindle book types are getting displayed in Product details page"
$ws.Range("B9").Value = "Verify that new book price and old book price (if available) are getting displayed in Paperback tab."
$ws.Range("B10").Value = "Verify that Kindle book details (Price, book features etc) are getting displayed after clicking on 'Kindle' tab in Product details page."
$ws.Range("B6").Value = "Verify that selected book details (Book title, authors, Book available types etc) are getting displayed in Product details page."
$ws.Range("B8").Value = "Verify that Paperback/Hardcover is displayed by default after navgating to Product details page."

# Give the new rows the same formatting (fill/border) as the existing
# data rows by copying the row-6 formats down.
$ws.Range("A6:B6").Copy()
$ws.Range("A8:B10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A8").Value = 7
$ws.Range("A9").Value = 8
$ws.Range("A10").Value = 9

# Widen column B to fit the longer text now in it. (Excel's ColumnWidth
# COM property only accepts/returns character widths quantized to 1/6ths,
# so 118.3 is the input that round-trips to the target stored width of
# ~119.14 characters.)
$ws.Columns.Item(2).ColumnWidth = 118.3
